# Certificate NR10 model: remove RG references from the "portador" sentence,
# shrink the textbox now that there is less text, and tidy up the
# "São Carlos, {{DATA}}" run split.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)   # "Rectangle 5" - the certificate body textbox
$tr = $shape.TextFrame.TextRange

# --- 1) Remove the RG mention from "portador do RG nº {{RG}} e CPF nº {{CPF}}, " ---
# Work right-to-left so earlier character offsets stay valid.

# "e CPF nº " (chars 53-61) -> "nº "
$rCpfNo = $tr.Characters(53, 9)
$rCpfNo.Text = "nº "

# "{{RG}} " (chars 46-52) -> "CPF "
$rRgPlaceholder = $tr.Characters(46, 7)
$rRgPlaceholder.Text = "CPF "

# tail of "portador do RG nº " (chars 37-45, i.e. "do RG nº ") -> "do "
$rDoRg = $tr.Characters(37, 9)
$rDoRg.Text = "do "

# --- 2) Merge "São Carlos" + ", " into a single run "São Carlos, " ---
# (offset shifted left by 15 chars because of the edits above)
$rSaoCarlos = $tr.Characters(243, 12)
$rSaoCarlos.Text = "São Carlos, "

# --- 3) Shrink the textbox height now that one line of text was removed ---
# (set after all text edits - the shape has spAutoFit, so resizing the
# text earlier in the script gets recomputed by autofit on later edits)
$shape.Height = 193.87503937007875
